$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1726.6428
$ws.Range("I2").Value = 1767.8572
$ws.Range("J2").Value = 1685.4286
$ws.Range("K2").Value = 1767.8572
$ws.Range("L2").Value = 1685.4286
$ws.Range("M2").Value = -1654.8572
$ws.Range("N2").Value = -1911.4286

# Row 19
$ws.Range("H19").Value = 997.8333
$ws.Range("I19").Value = 998.75
$ws.Range("J19").Value = 996
$ws.Range("K19").Value = 998.75
$ws.Range("L19").Value = 996
$ws.Range("M19").Value = -823.75
$ws.Range("N19").Value = -1346

# Row 86
$ws.Range("H86").Value = 6292.3887
$ws.Range("I86").Value = 3110.9
$ws.Range("K86").Value = 3110.9
$ws.Range("M86").Value = -1987.9

# Row 89
$ws.Range("H89").Value = 6292.3887
$ws.Range("I89").Value = 3110.9
$ws.Range("K89").Value = 15554.5
$ws.Range("M89").Value = -9938.5

# Row 98
$ws.Range("H98").Value = 1388.9333
$ws.Range("I98").Value = 1388.9333
$ws.Range("K98").Value = 1388.9333
$ws.Range("M98").Value = 109.0667000000001

# Row 106
$ws.Range("H106").Value = 281549.6
$ws.Range("I106").Value = 399806.94
$ws.Range("K106").Value = 399806.94
$ws.Range("M106").Value = -399175.94

# Row 113
$ws.Range("H113").Value = 9001.223
$ws.Range("I113").Value = 5005
$ws.Range("J113").Value = 9500.75
$ws.Range("K113").Value = 5005
$ws.Range("L113").Value = 9500.75
$ws.Range("M113").Value = -1751
$ws.Range("N113").Value = -16008.75

# Row 122
$ws.Range("H122").Value = 1388.9333
$ws.Range("I122").Value = 1388.9333
$ws.Range("K122").Value = 4166.7999
$ws.Range("M122").Value = -1716.7999

# Row 132
$ws.Range("H132").Value = 2248.3845
$ws.Range("I132").Value = 1111.7391
$ws.Range("J132").Value = 10962.667
$ws.Range("K132").Value = 3335.2173
$ws.Range("L132").Value = 32888.001
$ws.Range("M132").Value = -805.2173000000003
$ws.Range("N132").Value = -37948.001

# Row 133
$ws.Range("H133").Value = 92865.266
$ws.Range("J133").Value = 92865.266
$ws.Range("L133").Value = 92865.266
$ws.Range("N133").Value = -102985.266

# Row 138
$ws.Range("H138").Value = 2186.7625
$ws.Range("I138").Value = 1368.7693
$ws.Range("J138").Value = 2580.611
$ws.Range("K138").Value = 4106.3079
$ws.Range("L138").Value = 7741.833
$ws.Range("M138").Value = 1033.6921
$ws.Range("N138").Value = -18021.833

$ws = $wb.Worksheets.Item("ARM")
# Row 21
$ws.Range("H21").Value = 438
$ws.Range("I21").Value = 551.6667
$ws.Range("J21").Value = 267.5
$ws.Range("K21").Value = 551.6667
$ws.Range("L21").Value = 267.5
$ws.Range("M21").Value = -177.6667
$ws.Range("N21").Value = -1015.5

# Row 26
$ws.Range("H26").Value = 14127
$ws.Range("I26").Value = 12166.667
$ws.Range("K26").Value = 12166.667
$ws.Range("M26").Value = -11836.667

# Row 30
$ws.Range("H30").Value = 6693.143
$ws.Range("J30").Value = 6693.143
$ws.Range("L30").Value = 6693.143
$ws.Range("N30").Value = -6993.143

# Row 32
$ws.Range("H32").Value = 6801.33
$ws.Range("I32").Value = 3816.1204
$ws.Range("K32").Value = 3816.1204
$ws.Range("M32").Value = -3529.1204

# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

# Row 44
$ws.Range("H44").Value = 122288
$ws.Range("J44").Value = 122288
$ws.Range("L44").Value = 122288
$ws.Range("N44").Value = -123264

# Row 55
$ws.Range("H55").Value = 37142.855
$ws.Range("J55").Value = 58000
$ws.Range("L55").Value = 58000
$ws.Range("N55").Value = -58630

# Row 132
$ws.Range("H132").Value = 2968.75
$ws.Range("I132").Value = 2042.1333
$ws.Range("J132").Value = 4037.923
$ws.Range("K132").Value = 6126.3999
$ws.Range("L132").Value = 12113.769
$ws.Range("M132").Value = -3596.3999
$ws.Range("N132").Value = -17173.769

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 69964.60000000001
$ws.Range("J20").Value = 5067.2856
$ws.Range("L20").Value = 5067.2856
$ws.Range("N20").Value = -5561.2856

# Row 22
$ws.Range("H22").Value = 19231518
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 5562.5
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 7083.3335
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 7083.3335
$ws.Range("M15").Value = -830
$ws.Range("N15").Value = -7423.3335

# Row 31
$ws.Range("H31").Value = 5523
$ws.Range("I31").Value = 3361.5715
$ws.Range("K31").Value = 3361.5715
$ws.Range("M31").Value = -3066.5715

# Row 34
$ws.Range("H34").Value = 5523
$ws.Range("I34").Value = 3361.5715
$ws.Range("K34").Value = 3361.5715
$ws.Range("M34").Value = -3159.5715

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1981.5264
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 15000
$ws.Range("N5").Value = -15224

# Row 14
$ws.Range("H14").Value = 505.14285
$ws.Range("I14").Value = 505.14285
$ws.Range("K14").Value = 1515.42855
$ws.Range("M14").Value = -1342.42855

# Row 134
$ws.Range("H134").Value = 91647.55
$ws.Range("I134").Value = 91647.55
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 274942.65
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -269872.65
$ws.Range("N134").ClearContents()

# Row 135
$ws.Range("H135").Value = 1981.5264
$ws.Range("J135").Value = 5000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -50070

# Row 140
$ws.Range("H140").Value = 20835422
$ws.Range("J140").Value = 3999
$ws.Range("L140").Value = 11997
$ws.Range("N140").Value = -22357

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 69816.31
$ws.Range("I70").Value = 83542.92
$ws.Range("J70").Value = 10334.333
$ws.Range("K70").Value = 83542.92
$ws.Range("L70").Value = 10334.333
$ws.Range("M70").Value = -83272.92
$ws.Range("N70").Value = -10874.333

# Row 73
$ws.Range("H73").Value = 69816.31
$ws.Range("I73").Value = 83542.92
$ws.Range("J73").Value = 10334.333
$ws.Range("K73").Value = 83542.92
$ws.Range("L73").Value = 10334.333
$ws.Range("M73").Value = -82606.92
$ws.Range("N73").Value = -12206.333

# Row 102
$ws.Range("H102").Value = 989.0303
$ws.Range("I102").Value = 994.9375
$ws.Range("K102").Value = 994.9375
$ws.Range("M102").Value = 627.0625

# Row 132
$ws.Range("H132").Value = 1574.5
$ws.Range("I132").Value = 1155.2354
$ws.Range("K132").Value = 3465.7062
$ws.Range("M132").Value = -935.7062000000001

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# Row 140
$ws.Range("H140").Value = 74495

# Row 141
$ws.Range("H141").Value = 125900
$ws.Range("J141").Value = 164750
$ws.Range("L141").Value = 164750
$ws.Range("N141").Value = -175110

$ws = $wb.Worksheets.Item("LTW")
# Row 26
$ws.Range("H26").Value = 1903.3334
$ws.Range("J26").Value = 2355
$ws.Range("L26").Value = 2355
$ws.Range("N26").Value = -2945

# Row 130
$ws.Range("H130").Value = 65000
$ws.Range("J130").Value = 65000
$ws.Range("L130").Value = 65000
$ws.Range("N130").Value = -75040

# Row 131
$ws.Range("H131").Value = 42249
$ws.Range("J131").Value = 42249
$ws.Range("L131").Value = 42249
$ws.Range("N131").Value = -52329

# Row 132
$ws.Range("H132").Value = 3342.8823
$ws.Range("I132").Value = 2874.4285
$ws.Range("J132").Value = 3670.8
$ws.Range("K132").Value = 8623.2855
$ws.Range("L132").Value = 11012.4
$ws.Range("M132").Value = -6093.2855
$ws.Range("N132").Value = -16072.4

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 11244.954
$ws.Range("I126").Value = 2177.8333
$ws.Range("J126").Value = 22125.5
$ws.Range("K126").Value = 6533.499899999999
$ws.Range("L126").Value = 66376.5
$ws.Range("M126").Value = -4063.499899999999
$ws.Range("N126").Value = -71316.5
